$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (personnel_id 3) was "Alex Schmidt" -> now "Max Schmidt".
# Surname/street/zip/city/password/role/isLoggedIn stay the same; only
# forename, username and email change.
$ws.Range("C4").Value = "Max"
$ws.Range("B4").Value = "max_schmidt"
$ws.Range("H4").Value = "max.schmidt@example.com"

# Selection moved from N16 to M16.
$ws.Range("M16").Select() | Out-Null
